$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 2, shifting the existing rows (old 2,3,4 -> new 3,4,5)
# down. This makes room for a new "slug" row that lets two columns relate to each
# other via machine-friendly identifiers (fixes #13: hierarchical SKOS linking).
$ws.Rows.Item(2).Insert()

# Fill the new row with the "slug" identifiers for each column.
$ws.Range("A2").Value = "sector-descripcion"
$ws.Range("B2").Value = "sector-codigo"
$ws.Range("C2").Value = "n-parados"
$ws.Range("D2").Value = "comarca-nombre"
$ws.Range("E2").Value = "comarca-codigo"
$ws.Range("F2").Value = "sector-descripcion"
$ws.Range("G2").Value = "aragon"
$ws.Range("H2").Value = "n-demandantes"
$ws.Range("I2").Value = "provincia-codigo"
$ws.Range("J2").Value = "provincia-nombre"
$ws.Range("K2").Value = "sexo"
$ws.Range("L2").Value = "mes-y-ano"
